$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "33.987.85"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.781.99"
$ws.Range("E3").Value = "  +0.37%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "225.84"
$ws.Range("E5").Value = "  +2.19%  "

# Row 6 - XRP
Set-TextValue "D6" "0.553"
$ws.Range("E6").Value = "  +1.28%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8 - Solana
Set-TextValue "D8" "32.23"
$ws.Range("E8").Value = "  +3.17%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.291"
$ws.Range("E9").Value = "  +1.89%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0704"
$ws.Range("E10").Value = "  -0.01%  "

# Row 11 - TRON (price unchanged)
$ws.Range("E11").Value = "  +1.77%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.040.08"
$ws.Range("E12").Value = "  +0.52%  "

# Row 13 - was Chainlink, now WrappedEther (rows 13 & 14 swapped content)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.830.41"
$ws.Range("E13").Value = "  +2.88%  "

# Row 14 - was WrappedEther, now Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "11.05"
$ws.Range("E14").Value = "  +4.95%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.623"
$ws.Range("E15").Value = "  +0.36%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "33.955.62"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17 - Polkadot (price unchanged)
$ws.Range("E17").Value = "  -0.77%  "

# Row 18 - Litecoin
Set-TextValue "D18" "67.90"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "243.05"
$ws.Range("E19").Value = "  -0.05%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0784"
$ws.Range("E20").Value = "  +1.56%  "

# Row 21 - Dai
Set-TextValue "D21" "0.998"
$ws.Range("E21").Value = "  -0.25%  "

# Row 22 - Avalanche
Set-TextValue "D22" "10.72"
$ws.Range("E22").Value = "  +1.86%  "

# Row 23 - Uniswap (price unchanged)
$ws.Range("E23").Value = "  +1.68%  "

# Row 24 - Toncoin (price unchanged)
$ws.Range("E24").Value = "  -2.87%  "

# Row 25 - Monero
Set-TextValue "D25" "159.72"
$ws.Range("E25").Value = "  +1.05%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "16.29"
$ws.Range("E26").Value = "  -0.26%  "

# Row 27 - Cosmos
Set-TextValue "D27" "7.10"
$ws.Range("E27").Value = "  +1.54%  "

# Row 28 - Stellar (price unchanged)
$ws.Range("E28").Value = "  +1.49%  "

# Row 29 - BinanceUSD (price unchanged)
$ws.Range("E29").Value = "  +0.21%  "

# Row 30 - PancakeSwap (price unchanged)
$ws.Range("E30").Value = "  +3.71%  "

# Row 31 - Hedera
Set-TextValue "D31" "0.0512"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.66"
$ws.Range("E32").Value = "  -0.60%  "

# Row 33 - InternetComputer(DFINITY) (price unchanged)
$ws.Range("E33").Value = "  +1.40%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.81"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35 - Maker
Set-TextValue "D35" "1.393.53"
$ws.Range("E35").Value = "  -0.05%  "

# Row 36 - ImmutableX (price unchanged)
$ws.Range("E36").Value = "  +5.62%  "

# Row 37 - TrustWalletToken (price unchanged)
$ws.Range("E37").Value = "  -0.41%  "

# Row 38 - VeChain (price unchanged)
$ws.Range("E38").Value = "  +1.12%  "

# Row 39 - RenderToken
Set-TextValue "D39" "2.24"
$ws.Range("E39").Value = "  +7.03%  "

# Row 40 - HuobiToken (price unchanged)
$ws.Range("E40").Value = "  +1.07%  "

# Row 41 - ARBITRUM
Set-TextValue "D41" "0.914"
$ws.Range("E41").Value = "  -1.75%  "

# Row 42 - Aave
Set-TextValue "D42" "77.92"
$ws.Range("E42").Value = "  -1.21%  "

# Row 43 - MXToken
Set-TextValue "D43" "2.66"
$ws.Range("E43").Value = "  -1.09%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "13.48"
$ws.Range("E44").Value = "  +14.94%  "

# Row 45 - BabyDogeCoin (price unchanged)
$ws.Range("E45").Value = "  +19.24%  "

# Row 47 - Quant
Set-TextValue "D47" "108.33"
$ws.Range("E47").Value = "  +4.71%  "

# Row 48 - Kaspa (price unchanged)
$ws.Range("E48").Value = "  +1.20%  "

# Row 49 - FraxShare
Set-TextValue "D49" "5.85"
$ws.Range("E49").Value = "  +0.48%  "

# Row 50 - RocketPoolETH
Set-TextValue "D50" "1.939.59"
$ws.Range("E50").Value = "  +0.62%  "

# Row 51 - PaxDollar
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  +0.58%  "
